# This workbook's data rows got re-shuffled upstream: the content that used
# to live in one row now lives in another (two independent cycles: rows
# 2/3/23, and rows 10-24 excluding 23). Row 1 (header) and rows 4-9 are not
# part of either cycle and stay untouched.
#
# Strategy: snapshot every source row's full A:AY values first (so we never
# read a row after it has already been overwritten), then write each
# snapshot into its destination row.
#
# A handful of columns (Y, Z, AA, AB) hold date/time values stored as plain
# text (e.g. "2011-10-23", "00:00") rather than real Excel dates. Writing a
# date-looking string straight into a General-formatted cell makes Excel
# auto-coerce it into a date serial, which would change the cell's stored
# type. Pre-formatting those columns as Text ("@") keeps the round-tripped
# values as plain strings, matching the source representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 24
$ws.Range("Y1:Y$lastRow").NumberFormat = "@"
$ws.Range("Z1:Z$lastRow").NumberFormat = "@"
$ws.Range("AA1:AA$lastRow").NumberFormat = "@"
$ws.Range("AB1:AB$lastRow").NumberFormat = "@"

# newRow -> oldRow: the row whose content should end up at newRow.
$mapping = @{
  2  = 23
  3  = 2
  10 = 11
  11 = 12
  12 = 13
  13 = 14
  14 = 15
  15 = 16
  16 = 17
  17 = 18
  18 = 19
  19 = 20
  20 = 21
  21 = 22
  22 = 24
  23 = 3
  24 = 10
}

# Snapshot each distinct source row exactly once, before any writes happen.
$snapshots = @{}
foreach ($oldRow in $mapping.Values) {
    if (-not $snapshots.ContainsKey($oldRow)) {
        $snapshots[$oldRow] = $ws.Range("A" + $oldRow + ":AY" + $oldRow).Value()
    }
}

# Now write every snapshot to its new destination.
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $ws.Range("A" + $newRow + ":AY" + $newRow).Value = $snapshots[$oldRow]
}
